# "save game plus tekniske dokumenter"
# Moves the "MYSQL" info block up (from rows 16-19 to rows 14-17), renames its
# header to "MYSQL @", and adds the database server's IP address next to the
# new header - mirroring the "SERVER @" / "DjATABASE @" blocks above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "MYSQL" block entirely (header + the three detail rows below it).
$ws.Range("B16:C19").Clear()

# Re-create the block starting two rows higher, with the updated header text
# and a leading IP-address column (like the other two blocks on the sheet).
$ws.Range("B14").Value = "MYSQL @"
$ws.Range("B14").Font.Bold = $true
$ws.Range("C14").Value = "10.2.1.98"

$ws.Range("B15").Value = "Login MySQL"
$ws.Range("C15").Value = "sudo mysql -u root -p"

$ws.Range("B16").Value = "root -p"
$ws.Range("C16").Value = "r00tc4nrun"

$ws.Range("B17").Value = "phpMyAdmin"
$ws.Range("C17").Value = "pyh0ypy"

# Leave the selection where the author last left it before saving.
[void]$ws.Range("B18").Select()
